$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the Honduras Liga Nacional match data (previously on row 3),
# using the newly updated odds values.
$ws.Range("A2").Value = "Honduras Liga Nacional"

# B2/C2 hold date/time-looking text ("2026-01-07" / "22:00:00"). Excel would
# otherwise auto-convert these into real date/time serial values (and apply a
# number format), so temporarily mark the cells as Text, assign the value,
# then restore the default "Normal" style so no stray style is left behind.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2026-01-07"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "22:00:00"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "CD Olimpia"
$ws.Range("E2").Value = "CD Marathon"
$ws.Range("F2").Value = 1.9
$ws.Range("G2").Value = 1.96
$ws.Range("H2").Value = 4.9
$ws.Range("I2").Value = 5.4
$ws.Range("J2").Value = 3.35
$ws.Range("K2").Value = 3.55
$ws.Range("L2").Value = 1.46
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.37
$ws.Range("P2").Value = 1.84
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.31
$ws.Range("S2").Value = 3.85
$ws.Range("T2").Value = 1.93
$ws.Range("U2").Value = 1.89
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 2.04
$ws.Range("X2").Value = 12.5
$ws.Range("Y2").Value = 17.5
$ws.Range("Z2").Value = 38
$ws.Range("AA2").Value = 140
$ws.Range("AB2").Value = 8.4
$ws.Range("AC2").Value = 8.4
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 85
$ws.Range("AF2").Value = 11
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 120
$ws.Range("AJ2").Value = 23
$ws.Range("AK2").Value = 22
$ws.Range("AL2").Value = 60
$ws.Range("AM2").Value = 150
$ws.Range("AN2").Value = 17.5
$ws.Range("AO2").Value = 160

# Remove the now-duplicate row 3 (it held the same match before the update),
# shrinking the used range to A1:AO2.
$ws.Rows.Item(3).Delete()
